$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cinema")

# Insert a new column at S (shifts old S..AML right by one, to T..AMM).
# Excel's entire-column Insert copies formatting from the column to the
# left (R) onto the freshly inserted column, matching the target layout.
[void]$ws.Columns("S:S").Insert()

# Header for the newly inserted column.
$ws.Range("S1").Value = "Sub brand"

# The AutoFilter range must grow to cover the new column (A1:AO54 -> A1:AP54).
$ws.AutoFilterMode = $false
[void]$ws.Range("A1:AP54").AutoFilter()

# The workbook-level defined names tied to the filter database must be
# resized the same way.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Cinema!_FilterDatabase" -or $n.Name -eq "Cinema!_FilterDatabase_0") {
        $n.RefersTo = "=Cinema!`$A`$1:`$AP`$54"
    }
}

# Restore the active selection to S2, as recorded in the saved view state.
[void]$ws.Range("S2").Select()
